# Mark the following VK API methods as implemented on the "Methods" sheet:
#   groups.getBanned, groups.getById, groups.getMembers, groups.isMember,
#   messages.createChat, messages.delete
# Their Status column (C) moves from "In progress" to "Implemeted",
# matching the formatting already used by other "Implemeted" rows.

$wb = $excel.ActiveWorkbook
$wsMethods = $wb.Worksheets.Item("Methods")
$wsStats = $wb.Worksheets.Item("Statistics")

# Use an existing "Implemeted" cell as the formatting template.
$template = $wsMethods.Range("C8")
$templateStyle = $template.Style
$templateNumberFormat = $template.NumberFormat
$templateHAlign = $template.HorizontalAlignment
$templateVAlign = $template.VerticalAlignment

$rows = @(23, 24, 30, 33, 40, 41)
foreach ($r in $rows) {
    $cell = $wsMethods.Cells.Item($r, 3)
    $cell.Value = "Implemeted"
    $cell.Style = $templateStyle
    $cell.NumberFormat = $templateNumberFormat
    $cell.HorizontalAlignment = $templateHAlign
    $cell.VerticalAlignment = $templateVAlign
}

# Reflect the newly-finished methods in the workbook's recalculated totals
# (Statistics!B4 / B5 depend on these via COUNTIFS and recalc automatically).
$wsStats.Calculate()

# Make "Methods" the active sheet/tab with cell C41 selected, as in the
# edited workbook.
$wsMethods.Activate()
$wsMethods.Range("C41").Select()
